$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark row 10 ("Find Sum of AP Series") as done, matching the other rows in column B
$ws.Range("B10").Value = "D"

# Update the saved view state: select C13 and scroll so row 7 is the top-left row
$ws.Activate()
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
